$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.598.61"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.114.68"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.014"
$ws.Range("E4").Value = "  +1.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "342.94"
$ws.Range("E5").Value = "  +2.76%  "
$ws.Range("E6").Value = "  +1.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5251"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4505"
$ws.Range("E8").Value = "  -1.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.51"
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09016"
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.170"
$ws.Range("E11").Value = "  -1.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.35"
$ws.Range("E12").Value = "  -0.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.112.41"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.798"
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.066"
$ws.Range("E15").Value = "  +2.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "98.20"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001169"
$ws.Range("E17").Value = "  +3.27%  "
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06712"
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.34"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.011"
$ws.Range("E21").Value = "  +1.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.331"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.677.63"
$ws.Range("E23").Value = "  +0.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.76"
$ws.Range("E24").Value = "  +3.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.378"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.360.79"
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.30"
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.88"
$ws.Range("E28").Value = "  +1.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.533"
$ws.Range("E29").Value = "  -2.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.77"
$ws.Range("E30").Value = "  +1.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.194"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1072"
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.356"
$ws.Range("E33").Value = "  +3.25%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.630"
$ws.Range("E34").Value = "  -2.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.955"
$ws.Range("E35").Value = "  +0.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.25"
$ws.Range("E36").Value = "  -1.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.845"
$ws.Range("E37").Value = "  +4.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02644"
$ws.Range("E38").Value = "  +2.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06826"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2329"
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.64"
$ws.Range("E41").Value = "  -0.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6854"
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("E43").Value = "  +1.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.75"
$ws.Range("E44").Value = "  +4.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6409"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.312"
$ws.Range("E46").Value = "  -1.77%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.727"
$ws.Range("E47").Value = "  +2.00%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000367"
$ws.Range("E48").Value = "  +3.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.254"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.79"
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07293"
$ws.Range("E51").Value = "  +2.80%  "
